$p = $ppt.ActivePresentation

# Add a new slide at the end, using the Blank layout (same as the other content slides)
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 12)

# Add the 6-row x 5-column table shape
$shp = $s.Shapes.AddTable(6, 5, 339281, 241300, 11330559, 2397760)
$tbl = $shp.Table

$tbl.Columns.Item(1).Width = 997839
$tbl.Columns.Item(2).Width = 6400800
$tbl.Columns.Item(3).Width = 1188720
$tbl.Columns.Item(4).Width = 1097280
$tbl.Columns.Item(5).Width = 1645920

# Row 1 - header
$tbl.Cell(1,1).Shape.TextFrame.TextRange.Text = "API"
$tbl.Cell(1,2).Shape.TextFrame.TextRange.Text = "huntss.py"

# Row 2 - method
$tbl.Cell(2,1).Shape.TextFrame.TextRange.Text = "method"
$tbl.Cell(2,2).Shape.TextFrame.TextRange.Text = "get_all_active"

# Row 3 - route
$tbl.Cell(3,1).Shape.TextFrame.TextRange.Text = "route"
$tbl.Cell(3,2).Shape.TextFrame.TextRange.Text = "/hunts/active [GET]"

# Row 4 - frontend usage
$tbl.Cell(4,1).Shape.TextFrame.TextRange.Text = "frontend `r`nusage"
$tbl.Cell(4,2).Shape.TextFrame.TextRange.Text = "hunt"

# Row 5 - results
$tbl.Cell(5,1).Shape.TextFrame.TextRange.Text = "results"
$tbl.Cell(5,2).Shape.TextFrame.TextRange.Text = "list of all hunts that aren't closed (should only ever be 1)"
$tbl.Cell(5,3).Shape.TextFrame.TextRange.Text = "Redis key"
$tbl.Cell(5,3).Shape.TextFrame.TextRange.Font.Bold = 1
$tbl.Cell(5,3).Shape.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$tbl.Cell(5,4).Shape.TextFrame.TextRange.Text = "Expiration"
$tbl.Cell(5,4).Shape.TextFrame.TextRange.Font.Bold = 1
$tbl.Cell(5,4).Shape.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$tbl.Cell(5,5).Shape.TextFrame.TextRange.Text = "Delete"
$tbl.Cell(5,5).Shape.TextFrame.TextRange.Font.Bold = 1
$tbl.Cell(5,5).Shape.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# Row 6 - DB calls
$tbl.Cell(6,1).Shape.TextFrame.TextRange.Text = "DB calls"
$tbl.Cell(6,2).Shape.TextFrame.TextRange.Text = "SELECT * FROM hunts WHERE status != 'hunt_closed' ORDER BY hunt_date"
$tbl.Cell(6,3).Shape.TextFrame.TextRange.Text = "sierra"
$tbl.Cell(6,3).Shape.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$tbl.Cell(6,4).Shape.TextFrame.TextRange.Text = "1 day"
$tbl.Cell(6,4).Shape.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$tbl.Cell(6,5).Shape.TextFrame.TextRange.Text = "hunts CUD"
